$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$reqLOB = "LOB1053 -  Física III  (Requisito)`n"
$reqLOM = "LOM3202 -  Circuitos Elétricos  (Indicação de Conjunto)`n"

# Swap the two requirement strings: LOM3202 now comes first (row 24), LOB1053 second (row 25)
$ws.Range("B24").Value = $reqLOM
$ws.Range("C24").Value = $reqLOM
$ws.Range("B25").Value = $reqLOB
$ws.Range("C25").Value = $reqLOB
